# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45205 to 45206 (2023-10-06 -> 2023-10-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data rows run from row 2 to row 482 in column C.
$ws.Range("C2:C482").Value = 45206
